$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("archive_formats_by_aip_2021-08")

# Clear the NARA_PRONOM URL (M3) and NARA_Proposed Preservation Plan (O3) cells
# for row 3, per merge_format_reports.py output: these should be blank (NaN)
# instead of "NO VALUE".
$ws.Range("M3").ClearContents()
$ws.Range("O3").ClearContents()

# Select L1:P1 as the active selection, matching the saved sheet view state.
$ws.Range("L1:P1").Select()
